# Applies targeted odds/value updates to Sheet1 of the FlashScore workbook.
# Each assignment below corresponds to a single cell change captured in the
# authoritative diff (same sheet, "Sheet1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BD5").Value = 151

$ws.Range("G10").Value = 2.5
$ws.Range("I10").Value = 2.9
$ws.Range("Y10").Value = 10
$ws.Range("AA10").Value = 21
$ws.Range("AW10").Value = 4.75

$ws.Range("G17").Value = 2.7
$ws.Range("J17").Value = 3.25
$ws.Range("N17").Value = 7.8
$ws.Range("Q17").Value = 1.78
$ws.Range("R17").Value = 1.98
$ws.Range("S17").Value = 1.37
$ws.Range("T17").Value = 2.87
$ws.Range("W17").Value = 9.75
$ws.Range("Z17").Value = 30
$ws.Range("AA17").Value = 21
$ws.Range("AB17").Value = 27
$ws.Range("AC17").Value = 7.8
$ws.Range("AI17").Value = 12.5
$ws.Range("AL17").Value = 18.5
$ws.Range("AO17").Value = 14.5
$ws.Range("AP17").Value = 21
$ws.Range("AQ17").Value = 60
$ws.Range("AR17").Value = 90
$ws.Range("AT17").Value = 2.87
$ws.Range("AW17").Value = 4.45
$ws.Range("BA17").Value = 80

$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 3.8
$ws.Range("I28").Value = 3.4
$ws.Range("J28").Value = 2.6
$ws.Range("S28").Value = 1.3
$ws.Range("T28").Value = 3.4
$ws.Range("U28").Value = 1.53
$ws.Range("V28").Value = 2.38
$ws.Range("Y28").Value = 9
$ws.Range("Z28").Value = 19
$ws.Range("AC28").Value = 15
$ws.Range("AH28").Value = 13
$ws.Range("AI28").Value = 19
$ws.Range("AT28").Value = 3.4
$ws.Range("AX28").Value = 17

$ws.Range("N29").Value = 17
$ws.Range("Q29").Value = 1.5
$ws.Range("R29").Value = 2.5

$ws.Range("G31").Value = 4
$ws.Range("H31").Value = 3.9
$ws.Range("I31").Value = 1.8
$ws.Range("L31").Value = 2.3
$ws.Range("U31").Value = 1.5
$ws.Range("V31").Value = 2.5
$ws.Range("AE31").Value = 12
$ws.Range("AH31").Value = 11
$ws.Range("AI31").Value = 11
$ws.Range("AJ31").Value = 8.5
$ws.Range("AN31").Value = 6.5
$ws.Range("AR31").Value = 67
$ws.Range("AX31").Value = 9
$ws.Range("AZ31").Value = 26

$ws.Range("I45").Value = 1.87

$ws.Range("I47").Value = 1.69

$ws.Range("I49").Value = 2.7
$ws.Range("Y49").Value = 9.5
$ws.Range("AD49").Value = 6.5
$ws.Range("AK49").Value = 29
$ws.Range("AN49").Value = 4.5
$ws.Range("AZ49").Value = 51

$ws.Range("G82").Value = 4.75
$ws.Range("H82").Value = 4.1
$ws.Range("I82").Value = 1.55
$ws.Range("U82").Value = 1.83
$ws.Range("V82").Value = 1.83
$ws.Range("Y82").Value = 17
$ws.Range("AD82").Value = 8.5
$ws.Range("AG82").Value = 700
$ws.Range("AM82").Value = 26
$ws.Range("AP82").Value = 34
$ws.Range("AU82").Value = 8.5
$ws.Range("AW82").Value = 3.6
